# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G16").Value = 1181084
$ws.Range("G17").Value = 1141144
$ws.Range("G18").Value = 1181084
$ws.Range("G19").Value = 1141144
$ws.Range("G20").Value = 1141144
